$d = $word.ActiveDocument

# ---- Sanity check ----
$origCount = $d.Paragraphs.Count
Write-Host "Starting paragraph count: $origCount"

# ==================================================================
# PASS 1: in-place text updates (paragraph count unchanged, so every
# index below refers to the ORIGINAL (before-edit) paragraph layout)
# ==================================================================
$d.Paragraphs.Item(4).Range.Text = "Distinguished Research & Data Analytics Professional with 21 years of expertise in survey methodology, consumer insights, voting behavior, and advanced data analysis. Proven track record in designing and implementing comprehensive research studies, managing cross-functional teams, and translating complex data into actionable intelligence. Expert in geospatial analysis, demographic segmentation, and consumer behavior modeling with experience serving major brands, organizations, and political candidates. Regular expert testimony and source on public opinion for journalists, with redistricting analysis used in court cases."
$d.Paragraphs.Item(6).Range.Text = "Survey Methodology & Research Design: Survey Design and Questionnaire Development for Political and Market Research • Sampling Methodology and Statistical Analysis (R, SPSS, Stata, OSCAR) • Random Device Engagement (RDE), Text Message, Web Panel, and Live Telephone Calling • Focus Groups and Qualitative Research Methodologies • Meta-analytical Dataset Development for Longitudinal Analysis • Survey Instrument Standardization and Call Methods Optimization • Expert Testimony and Consultation on Research Methodology"
$d.Paragraphs.Item(7).Range.Text = "Data Analysis & Visualization: Advanced Statistical Modeling and Analysis (Regression, Clustering, Segmentation) • Data Visualization: Tableau, PowerBI, Seaborn, Matplotlib, d3.js • Geospatial Analysis: ArcGIS, Quantum GIS, GRASS, OSGeo, PostGIS • Choropleths and Hexagonal Grid Maps for Demographic Visualization • Consumer Behavior Analysis and Market Segmentation • Machine Learning and Predictive Modeling for Targeting • Big Data Analytics: Spark/PySpark, Hadoop, Snowflake, dbt"
$d.Paragraphs.Item(8).Range.Text = "Research Leadership & Client Management: Multi-million Dollar Research Project Management • Cross-functional Team Leadership (Teams of 7-11 professionals) • Client Relationship Management across Political, NGO, and Corporate Sectors • Stakeholder Briefing for Elected Officials and Senior Leadership • Court Case Analysis and Expert Testimony • Research Framework Development and Quality Control • Business Intelligence and Market Intelligence Delivery"
$d.Paragraphs.Item(12).Range.Text = "• Conducted comprehensive quantitative and qualitative research studies for political candidates and major organizations, providing actionable consumer insights and market intelligence"
$d.Paragraphs.Item(13).Range.Text = "• Designed and implemented advanced segmentation models using demographic, psychographic, and behavioral data to identify high-value targets"
$d.Paragraphs.Item(14).Range.Text = "• Led multi-million dollar market research projects involving sensitive consumer data, ensuring compliance with privacy regulations while delivering actionable insights"
$d.Paragraphs.Item(15).Range.Text = "• Developed and deployed custom research software that processed billions of consumer records for pattern analysis, fraud detection and entity resolution"
$d.Paragraphs.Item(16).Range.Text = "• Built and maintained client relationships across diverse industries, consistently delivering insights that drove strategic decision-making"
$d.Paragraphs.Item(17).Range.Text = "• Architected and engineered cloud-based data warehouse solutions processing billions of records for electoral analytics and geospatial analysis"
$d.Paragraphs.Item(20).Range.Text = "• Conceived and led implementation of comprehensive multi-tenant data warehouse integrating consumer demographic, economic, and behavioral data"
$d.Paragraphs.Item(21).Range.Text = "• Overhauled the organization's survey methodology and polling operations, significantly improving data accuracy and response rates"
$d.Paragraphs.Item(22).Range.Text = "• Managed and developed one of the in-house polling teams, focusing on Random Device Engagement (RDE), text message and web panel collected surveys, with live telephone calling and focus groups"
$d.Paragraphs.Item(23).Range.Text = "• Worked on standardizing questions, survey instruments and call methods, along with building a meta-analytical dataset for longitudinal analysis"
$d.Paragraphs.Item(24).Range.Text = "• Managed a cross-functional team of eleven data engineers and analysts, establishing best practices for research methodology and data analysis"
$d.Paragraphs.Item(33).Range.Text = "• Designed comprehensive survey instruments for specialized voting segments and niche markets"
$d.Paragraphs.Item(34).Range.Text = "• Developed sophisticated analytical products and reports that delivered actionable insights to clients"
$d.Paragraphs.Item(35).Range.Text = "• Co-developed RACSO web application to manage all aspects of survey operations, from instrument design to data collection and analysis"
$d.Paragraphs.Item(36).Range.Text = "• Wrote RFP and analyzed bids from 1,200 vendors before selecting implementation partner for RACSO platform"
$d.Paragraphs.Item(39).Range.Text = "• Engineered FLEEM web application using Twilio's API to make thousands of simultaneous phone calls for IVR polls"
$d.Paragraphs.Item(40).Range.Text = "• Used FLEEM for early quantitative research in support of Senators Martin Heinrich and Elizabeth Warren"
$d.Paragraphs.Item(41).Range.Text = "• Led all aspects of survey design, implementation, data analysis, and reporting for major national studies"
$d.Paragraphs.Item(42).Range.Text = "• Developed new statistical methods for boundary estimation techniques, enhancing geographic market segmentation capabilities"
$d.Paragraphs.Item(45).Range.Text = "• Developed software solutions for political campaigns and advocacy groups using modern web technologies"
$d.Paragraphs.Item(46).Range.Text = "• Built web applications for voter engagement and campaign management with real-time data integration"
$d.Paragraphs.Item(47).Range.Text = "• Integrated third-party APIs and data sources for campaign tools and voter database management"
$d.Paragraphs.Item(48).Range.Text = "• Collaborated with political strategists to translate requirements into technical solutions using agile methodologies"
$d.Paragraphs.Item(51).Range.Text = "• Integrated technology solutions within organizational frameworks for social justice organizations using open source technologies"
$d.Paragraphs.Item(52).Range.Text = "• Developed data management systems for community organizing efforts with focus on accessibility and usability"
$d.Paragraphs.Item(53).Range.Text = "• Provided technical training and support to nonprofit staff on database management and data analysis tools"
$d.Paragraphs.Item(54).Range.Text = "• Built custom applications for community engagement and advocacy using web technologies and mobile platforms"
$d.Paragraphs.Item(57).Range.Text = "• Worked on all aspects of questionnaire design, sampling, reporting and analysis for political actors in Congressional, Senate and Presidential elections"
$d.Paragraphs.Item(58).Range.Text = "• Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party and affiliated actors"
$d.Paragraphs.Item(59).Range.Text = "• Developed polling consortium database that later became the Polling Consortium Database at The Analyst Institute"
$d.Paragraphs.Item(60).Range.Text = "• Designed questionnaires and analyzed data for complex market research studies across diverse industries"
$d.Paragraphs.Item(63).Range.Text = "• Administered all quantitative and qualitative research, ensuring that reporting was accurate and comprehensive"
$d.Paragraphs.Item(64).Range.Text = "• Managed all aspects of survey fielding for a multi-million dollar research firm, including scheduling, oversight, sampling, and quality control"
$d.Paragraphs.Item(65).Range.Text = "• Developed and implemented data warehousing solutions for efficient storage and retrieval of research findings"
$d.Paragraphs.Item(66).Range.Text = "• Created custom reports and data visualizations based on specific client requirements"
$d.Paragraphs.Item(68).Range.Text = "Survey Methodology & Research Innovation"
$d.Paragraphs.Item(69).Range.Text = "• Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party, later becoming the Polling Consortium Database at The Analyst Institute"
$d.Paragraphs.Item(70).Range.Text = "• Developed RACSO platform for pollsters to fully administer research, analyzing bids from 1,200 vendors before selecting implementation partner"
$d.Paragraphs.Item(71).Range.Text = "• Engineered FLEEM system using Twilio API for thousands of simultaneous phone calls for IVR polls supporting Senators Martin Heinrich and Elizabeth Warren"
$d.Paragraphs.Item(72).Range.Text = "• Pioneered the integration of advanced mapping techniques into standard reports, including choropleths and hexagonal grid maps"

# ==================================================================
# PASS 2: new paragraph insertions. Anchors refer to the ORIGINAL
# (before-edit) paragraph indices, which are still valid positions
# here because we always insert right after a fixed anchor going
# from the bottom of the document upward, and -- for several new
# paragraphs sharing one anchor -- in reverse (last-wanted-first)
# order, so the final on-page order reads top-to-bottom correctly.
# ==================================================================
# --- inserting 10 paragraph(s) after original paragraph 72 ---
$anchor = $d.Paragraphs.Item(72)
$anchor.Range.InsertParagraphAfter()
$d.Paragraphs.Item(72 + 1).Range.Text = "• Developed advanced data pipelines for machine learning applications enhancing consumer segmentation and predictive modeling"
$anchor = $d.Paragraphs.Item(72)
$anchor.Range.InsertParagraphAfter()
$d.Paragraphs.Item(72 + 1).Range.Text = "• Led multi-million dollar market research projects involving sensitive consumer data with privacy compliance"
$anchor = $d.Paragraphs.Item(72)
$anchor.Range.InsertParagraphAfter()
$d.Paragraphs.Item(72 + 1).Range.Text = "• Designed, architected and created multi-tenant data warehouse tracking decades of political, geographical, econometric change"
$anchor = $d.Paragraphs.Item(72)
$anchor.Range.InsertParagraphAfter()
$d.Paragraphs.Item(72 + 1).Range.Text = "• Conceived, architected, engineered and deployed cloud-based redistricting software used by thousands of analysts nationwide"
$styleSrc = $null
for ($si = 1; $si -le $d.Paragraphs.Count; $si++) {
    if ($d.Paragraphs.Item($si).Style.NameLocal -eq "Heading 3") {
        $styleSrc = $d.Paragraphs.Item($si)
        break
    }
}
$styleSrc.Range.Copy()
$anchor = $d.Paragraphs.Item(72)
$anchor.Range.InsertParagraphAfter()
$blank = $d.Paragraphs.Item(72 + 1)
$pasteSpot = $d.Range($blank.Range.Start, $blank.Range.Start)
$pasteSpot.Paste()
$trailingBlank = $d.Paragraphs.Item(72 + 2)
$trailingBlank.Range.Delete()
$d.Paragraphs.Item(72 + 1).Range.Text = "Data Infrastructure & Analytics"
$anchor = $d.Paragraphs.Item(72)
$anchor.Range.InsertParagraphAfter()
$d.Paragraphs.Item(72 + 1).Range.Text = "• Provided expert consultation on research methodology for diverse stakeholder groups"
$anchor = $d.Paragraphs.Item(72)
$anchor.Range.InsertParagraphAfter()
$d.Paragraphs.Item(72 + 1).Range.Text = "• Research analysis used in court cases addressing housing, redistricting, and community development"
$anchor = $d.Paragraphs.Item(72)
$anchor.Range.InsertParagraphAfter()
$d.Paragraphs.Item(72 + 1).Range.Text = "• Redistricting analysis used in court cases with rigorous methodology and expert testimony"
$anchor = $d.Paragraphs.Item(72)
$anchor.Range.InsertParagraphAfter()
$d.Paragraphs.Item(72 + 1).Range.Text = "• Regular expert testimony and source on public opinion for journalists, elected officials, and NGO leadership"
$styleSrc = $null
for ($si = 1; $si -le $d.Paragraphs.Count; $si++) {
    if ($d.Paragraphs.Item($si).Style.NameLocal -eq "Heading 3") {
        $styleSrc = $d.Paragraphs.Item($si)
        break
    }
}
$styleSrc.Range.Copy()
$anchor = $d.Paragraphs.Item(72)
$anchor.Range.InsertParagraphAfter()
$blank = $d.Paragraphs.Item(72 + 1)
$pasteSpot = $d.Range($blank.Range.Start, $blank.Range.Start)
$pasteSpot.Paste()
$trailingBlank = $d.Paragraphs.Item(72 + 2)
$trailingBlank.Range.Delete()
$d.Paragraphs.Item(72 + 1).Range.Text = "Expert Testimony & Court Cases"

# --- inserting 2 paragraph(s) after original paragraph 66 ---
$anchor = $d.Paragraphs.Item(66)
$anchor.Range.InsertParagraphAfter()
$d.Paragraphs.Item(66 + 1).Range.Text = "• Trained field staff on data collection protocols and quality control using standardized methodologies and best practices"
$anchor = $d.Paragraphs.Item(66)
$anchor.Range.InsertParagraphAfter()
$d.Paragraphs.Item(66 + 1).Range.Text = "• Introduced mapping and geospatial analysis into standard reporting procedures, enhancing the value of research deliverables"

# --- inserting 2 paragraph(s) after original paragraph 60 ---
$anchor = $d.Paragraphs.Item(60)
$anchor.Range.InsertParagraphAfter()
$d.Paragraphs.Item(60 + 1).Range.Text = "• Pioneered the integration of advanced mapping techniques into standard reports, including choropleths and hexagonal grid maps"
$anchor = $d.Paragraphs.Item(60)
$anchor.Range.InsertParagraphAfter()
$d.Paragraphs.Item(60 + 1).Range.Text = "• Conducted statistical modeling and analysis to address multifaceted consumer behavior questions"

# --- inserting 2 paragraph(s) after original paragraph 42 ---
$anchor = $d.Paragraphs.Item(42)
$anchor.Range.InsertParagraphAfter()
$d.Paragraphs.Item(42 + 1).Range.Text = "• Provided tabular and graphical reporting with plans for interactive data exploration capabilities"
$anchor = $d.Paragraphs.Item(42)
$anchor.Range.InsertParagraphAfter()
$d.Paragraphs.Item(42 + 1).Range.Text = "• Created comprehensive data visualization solutions that improved clients' understanding of complex research findings"

# --- inserting 2 paragraph(s) after original paragraph 36 ---
$anchor = $d.Paragraphs.Item(36)
$anchor.Range.InsertParagraphAfter()
$d.Paragraphs.Item(36 + 1).Range.Text = "• Standardized reporting methodologies to improve clarity and impact of research findings"
$anchor = $d.Paragraphs.Item(36)
$anchor.Range.InsertParagraphAfter()
$d.Paragraphs.Item(36 + 1).Range.Text = "• Introduced geospatial techniques to enhance market segmentation capabilities, providing clients with location-based consumer insights"

# --- inserting 1 paragraph(s) after original paragraph 24 ---
$anchor = $d.Paragraphs.Item(24)
$anchor.Range.InsertParagraphAfter()
$d.Paragraphs.Item(24 + 1).Range.Text = "• Developed advanced data pipelines for machine learning applications that enhanced consumer segmentation and predictive modeling capabilities"

# ==================================================================
# Final sanity check
# ==================================================================
$finalCount = $d.Paragraphs.Count
Write-Host "Final paragraph count: $finalCount"
Write-Host "Expected paragraph count: 91"
